# Update the date line and the multiplication problems throughout the
# document (table of two-digit-by-two-digit multiplication exercises).

$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-12 Wednesday", "2024-06-13 Thursday"),
    @("14×68=", "83×39="),
    @("73×46=", "13×79="),
    @("68×59=", "47×95="),
    @("25×58=", "83×56="),
    @("88×91=", "65×15="),
    @("27×34=", "76×23="),
    @("63×80=", "22×62="),
    @("28×32=", "61×86="),
    @("68×17=", "92×33="),
    @("41×61=", "23×47="),
    @("69×45=", "16×61="),
    @("61×77=", "91×87="),
    @("47×37=", "15×31="),
    @("32×11=", "82×42="),
    @("11×48=", "93×12="),
    @("52×47=", "23×34="),
    @("92×66=", "31×32="),
    @("87×57=", "73×73="),
    @("65×80=", "44×70="),
    @("47×99=", "38×74="),
    @("36×27=", "35×41="),
    @("32×62=", "19×53="),
    @("47×90=", "86×51="),
    @("37×77=", "43×63="),
    @("52×39=", "76×54=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
